$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "11×12=" "29×23="
Replace-Text "30×26=" "71×31="
Replace-Text "41×35=" "81×70="
Replace-Text "62×56=" "91×50="
Replace-Text "20×46=" "62×21="
Replace-Text "35×34=" "52×62="
Replace-Text "63×18=" "93×16="
Replace-Text "69×26=" "20×32="
Replace-Text "57×67=" "79×66="
Replace-Text "37×69=" "26×87="
Replace-Text "68×30=" "81×42="
Replace-Text "72×68=" "17×80="
Replace-Text "47×96=" "68×39="
Replace-Text "33×21=" "78×90="
Replace-Text "42×94=" "58×21="
Replace-Text "57×64=" "24×13="
Replace-Text "95×78=" "26×65="
Replace-Text "65×54=" "13×20="
Replace-Text "14×15=" "39×75="
Replace-Text "88×95=" "48×75="
Replace-Text "42×53=" "69×23="
Replace-Text "98×67=" "15×74="
Replace-Text "74×84=" "21×84="
Replace-Text "34×63=" "25×44="
Replace-Text "32×50=" "68×76="
